# Apply edits described in the diff:
# - Remove spaces around the "-" separator in the "具体时间范围" (E column) time-range strings
# - Update the "想去人数" (F column) numeric counts on sheets 展览 (1), 演出 (2), and 全部类型 (4)

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("E2").Value = "2024.03.16 10:30-03.16 17:00"
$ws1.Range("F2").Value = 101
$ws1.Range("E3").Value = "2024.03.16 09:00-03.17 17:00"
$ws1.Range("F3").Value = 2129
$ws1.Range("E4").Value = "2024.03.30 09:00-03.31 17:30"
$ws1.Range("F4").Value = 875
$ws1.Range("E5").Value = "2024.05.01 09:30-05.02 17:30"
$ws1.Range("F5").Value = 1386
$ws1.Range("E6").Value = "2024.06.09 10:00-06.10 17:00"
$ws1.Range("F6").Value = 371

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("E2").Value = "2024.03.29 20:00-03.29 21:30"
$ws2.Range("E3").Value = "2024.03.30 20:00-03.30 21:30"

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("E2").Value = "2024.03.16 10:30-03.16 17:00"
$ws4.Range("F2").Value = 101
$ws4.Range("E3").Value = "2024.03.16 09:00-03.17 17:00"
$ws4.Range("F3").Value = 2129
$ws4.Range("E4").Value = "2024.03.29 20:00-03.29 21:30"
$ws4.Range("E5").Value = "2024.03.30 20:00-03.30 21:30"
$ws4.Range("E6").Value = "2024.03.30 09:00-03.31 17:30"
$ws4.Range("F6").Value = 875
$ws4.Range("E7").Value = "2024.05.01 09:30-05.02 17:30"
$ws4.Range("F7").Value = 1386
$ws4.Range("E8").Value = "2024.06.09 10:00-06.10 17:00"
$ws4.Range("F8").Value = 371
